$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit rotates the data of rows 2, 4, 5 and 6 among themselves:
#   new row 2 <- old row 4
#   new row 4 <- old row 6
#   new row 5 <- old row 2
#   new row 6 <- old row 5
# Row 3 and row 7 (and the header row 1) are left untouched.

# ---------------------------------------------------------------
# Row 2 (becomes old row 4's data)
# ---------------------------------------------------------------
$ws.Range("A2").Value = 111835718
$ws.Range("B2").Value = 56398
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("M2").Value = "äldre spår"
$ws.Range("Q2").Value = 471101.0270993827
$ws.Range("R2").Value = 6810411.753755242
$ws.Range("S2").Value = 10

# ---------------------------------------------------------------
# Row 4 (becomes old row 6's data)
# ---------------------------------------------------------------
$ws.Range("A4").Value = 111835826
$ws.Range("Q4").Value = 470915.776864712
$ws.Range("R4").Value = 6810385.536630718
$ws.Range("S4").Value = 5
$ws.Range("AC4").Value = "även hackspettbo, troligen av tret"

# ---------------------------------------------------------------
# Row 5 (becomes old row 2's data)
# ---------------------------------------------------------------
$ws.Range("A5").Value = 111835745
$ws.Range("B5").Value = 77515
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 471152.5480076601
$ws.Range("R5").Value = 6810381.652036018

# ---------------------------------------------------------------
# Row 6 (becomes old row 5's data)
# ---------------------------------------------------------------
$ws.Range("A6").Value = 111835838
$ws.Range("B6").Value = 89423
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma"
$ws.Range("H6").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 470914.6782613794
$ws.Range("R6").Value = 6810368.79402096
$ws.Range("AC6").ClearContents()
$ws.Range("M6").ClearContents()
